# FichaTecnicaObras_1 - "reset password, mas estilos"
# Updates the text fields on slide 1 (Ficha Tecnica de Obras) to the new
# control number, date, location, electoral district, official and
# press/media data.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

function Set-ShapeText {
    param(
        $Shape,
        [string]$NewText
    )
    $tr = $Shape.TextFrame.TextRange
    # Replace the run text through a Characters() sub-range instead of
    # TextRange.Text directly: it edits the existing <a:r> in place and
    # keeps it free of an explicit <a:rPr>, matching the source runs that
    # rely on the paragraph's <a:defRPr> for formatting.
    $chars = $tr.Characters(1, $tr.Length)
    $chars.Text = $NewText
}

Set-ShapeText $s.Shapes.Item(4)  "SEGOB_002"
Set-ShapeText $s.Shapes.Item(6)  "2015-10-15"
Set-ShapeText $s.Shapes.Item(8)  "AGUASCALIENTES"
Set-ShapeText $s.Shapes.Item(9)  "Asientos"
Set-ShapeText $s.Shapes.Item(10) "I Distrito Electoral Federal de Aguascalientes"
Set-ShapeText $s.Shapes.Item(11) "ALEX"
Set-ShapeText $s.Shapes.Item(12) "PRI"
Set-ShapeText $s.Shapes.Item(13) "PRESINDENTE"

Set-ShapeText $s.Shapes.Item(16) "Arturo Gasca"
Set-ShapeText $s.Shapes.Item(17) "Secretario Municipal"
Set-ShapeText $s.Shapes.Item(18) "la falta de motivación de la población ha bajado la preferencia al partido"
Set-ShapeText $s.Shapes.Item(19) "Periódico"
Set-ShapeText $s.Shapes.Item(20) "Columna Política"
Set-ShapeText $s.Shapes.Item(21) "local"
Set-ShapeText $s.Shapes.Item(22) "3"
